# Scheduled-runner price/profit refresh across the Asura_Profits workbook.
# Updates currentAveragePrice* / Leve*Price* / LeveProfit* columns (H:N)
# for a set of affected leve rows on each crafting-job sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 261.7143
$ws.Range("I33").Value = 247.2
$ws.Range("J33").Value = 298
$ws.Range("K33").Value = 247.2
$ws.Range("L33").Value = 298
$ws.Range("M33").Value = -18.19999999999999
$ws.Range("N33").Value = -756
$ws.Range("H43").Value = 1046.1333
$ws.Range("I43").Value = 671.4286
$ws.Range("J43").Value = 1374
$ws.Range("K43").Value = 671.4286
$ws.Range("L43").Value = 1374
$ws.Range("M43").Value = -602.4286
$ws.Range("N43").Value = -1512
$ws.Range("H74").Value = 4155.5
$ws.Range("I74").Value = 4155.5
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 4155.5
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -3219.5
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 4155.5
$ws.Range("I77").Value = 4155.5
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 20777.5
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -16097.5
$ws.Range("N77").ClearContents()
$ws.Range("H113").Value = 2333.3333
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 1000
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 2254
$ws.Range("N113").Value = -9508
$ws.Range("H132").Value = 2586.8
$ws.Range("I132").Value = 1709.6666
$ws.Range("K132").Value = 5128.9998
$ws.Range("M132").Value = -2598.9998
$ws.Range("H138").Value = 2868.6897
$ws.Range("I138").Value = 1340.2273
$ws.Range("J138").Value = 3386.0154
$ws.Range("K138").Value = 4020.6819
$ws.Range("L138").Value = 10158.0462
$ws.Range("M138").Value = 1119.3181
$ws.Range("N138").Value = -20438.0462
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 17440
$ws.Range("J23").Value = 12142.857
$ws.Range("L23").Value = 12142.857
$ws.Range("N23").Value = -12660.857
$ws.Range("H32").Value = 11674.87
$ws.Range("I32").Value = 11674.87
$ws.Range("K32").Value = 11674.87
$ws.Range("M32").Value = -11387.87
$ws.Range("H37").Value = 22793
$ws.Range("I37").Value = 1134
$ws.Range("J37").Value = 30012.666
$ws.Range("K37").Value = 1134
$ws.Range("L37").Value = 30012.666
$ws.Range("M37").Value = -861
$ws.Range("N37").Value = -30558.666
$ws.Range("H44").Value = 39775
$ws.Range("J44").Value = 39775
$ws.Range("L44").Value = 39775
$ws.Range("N44").Value = -40751
$ws.Range("H63").Value = 5486.5557
$ws.Range("I63").Value = 5486.5557
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 5486.5557
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -4800.5557
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 5486.5557
$ws.Range("I66").Value = 5486.5557
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 27432.7785
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -24000.7785
$ws.Range("N66").ClearContents()
$ws.Range("H80").Value = 36000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 36000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 36000
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -37996
$ws.Range("H83").Value = 36000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 36000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 108000
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -117984
$ws.Range("H111").Value = 267322
$ws.Range("J111").Value = 267322
$ws.Range("L111").Value = 267322
$ws.Range("N111").Value = -275502
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 78514.14
$ws.Range("I82").Value = 78514.14
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 78514.14
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -78131.14
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 78514.14
$ws.Range("I85").Value = 78514.14
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 78514.14
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -77188.14
$ws.Range("N85").ClearContents()
$ws.Range("H99").Value = 1797.8572
$ws.Range("I99").Value = 1665.8334
$ws.Range("J99").Value = 2590
$ws.Range("K99").Value = 1665.8334
$ws.Range("L99").Value = 2590
$ws.Range("M99").Value = -167.8334
$ws.Range("N99").Value = -5586
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1798
$ws.Range("I16").Value = 1830.1666
$ws.Range("J16").Value = 1733.6666
$ws.Range("K16").Value = 1830.1666
$ws.Range("L16").Value = 1733.6666
$ws.Range("M16").Value = -1543.1666
$ws.Range("N16").Value = -2307.6666
$ws.Range("H113").Value = 1798
$ws.Range("I113").Value = 1830.1666
$ws.Range("J113").Value = 1733.6666
$ws.Range("K113").Value = 1830.1666
$ws.Range("L113").Value = 1733.6666
$ws.Range("M113").Value = 339.8334
$ws.Range("N113").Value = -6073.6666
$ws.Range("H132").Value = 1753.2727
$ws.Range("I132").Value = 1326.8572
$ws.Range("J132").Value = 2499.5
$ws.Range("K132").Value = 3980.5716
$ws.Range("L132").Value = 7498.5
$ws.Range("M132").Value = -1450.5716
$ws.Range("N132").Value = -12558.5
$ws.Range("H134").Value = 1992.4584
$ws.Range("I134").Value = 1683.8889
$ws.Range("K134").Value = 5051.6667
$ws.Range("M134").Value = -2516.6667
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 657.3333
$ws.Range("I2").Value = 823.53845
$ws.Range("J2").Value = 225.2
$ws.Range("K2").Value = 4941.2307
$ws.Range("L2").Value = 1351.2
$ws.Range("M2").Value = -4828.2307
$ws.Range("N2").Value = -1577.2
$ws.Range("H5").Value = 1423.6786
$ws.Range("J5").Value = 806
$ws.Range("L5").Value = 2418
$ws.Range("N5").Value = -2642
$ws.Range("H135").Value = 1423.6786
$ws.Range("J135").Value = 806
$ws.Range("L135").Value = 7254
$ws.Range("N135").Value = -12324
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1564.5294
$ws.Range("I113").Value = 1017.2
$ws.Range("J113").Value = 2346.4285
$ws.Range("K113").Value = 1017.2
$ws.Range("L113").Value = 2346.4285
$ws.Range("M113").Value = 1152.8
$ws.Range("N113").Value = -6686.4285
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1052.2727
$ws.Range("I22").Value = 445.83334
$ws.Range("J22").Value = 1780
$ws.Range("K22").Value = 445.83334
$ws.Range("L22").Value = 1780
$ws.Range("M22").Value = -150.83334
$ws.Range("N22").Value = -2370
$ws.Range("H27").Value = 1052.2727
$ws.Range("I27").Value = 445.83334
$ws.Range("J27").Value = 1780
$ws.Range("K27").Value = 445.83334
$ws.Range("L27").Value = 1780
$ws.Range("M27").Value = -338.83334
$ws.Range("N27").Value = -1994
$ws.Range("H93").Value = 1862.5
$ws.Range("J93").Value = 1883.3334
$ws.Range("L93").Value = 1883.3334
$ws.Range("N93").Value = -4379.3334
$ws.Range("H100").Value = 5244.6113
$ws.Range("I100").Value = 6854.8184
$ws.Range("K100").Value = 6854.8184
$ws.Range("M100").Value = -6313.8184
$ws.Range("H136").Value = 3098.6
$ws.Range("I136").Value = 2563.1738
$ws.Range("K136").Value = 7689.5214
$ws.Range("M136").Value = -5139.5214
